# Applies the "adds connection to change in meat demand for 3 regions" edit:
#  - Rename "EU27.C_yield_variation_PCT[Maize]"   -> "EU 27.C_yield_variation_PCT[Maize]"
#  - Rename "EU27.C_yield_variation_PCT[OilCrop]" -> "EU 27.C_yield_variation_PCT[OilCrop]"
#  - Add three new labeled rows (21, 23, 25) with value 0 in column B:
#       Global.AP demand mod pct[USA 1]
#       Global.AP demand mod pct[EU27 1]
#       Global.AP demand mod pct[CHIHKG 1]

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old "EU27..." labels first (they get re-entered below with the
# corrected "EU 27" spelling once the shared-string slot has been freed).
$ws.Range("A11").ClearContents()
$ws.Range("A13").ClearContents()

# New rows for the meat-demand connection
$ws.Range("A21").Value = "Global.AP demand mod pct[USA 1]"
$ws.Range("B21").Value = 0

$ws.Range("A23").Value = "Global.AP demand mod pct[EU27 1]"
$ws.Range("B23").Value = 0

$ws.Range("A25").Value = "Global.AP demand mod pct[CHIHKG 1]"
$ws.Range("B25").Value = 0

# Fix the EU27 labels to "EU 27" (space added)
$ws.Range("A11").Value = "EU 27.C_yield_variation_PCT[Maize]"
$ws.Range("A13").Value = "EU 27.C_yield_variation_PCT[OilCrop]"

# Match the selection noted in the saved file
$ws.Range("I25").Select()
